# Applies the "cryptos list" refresh for Thu Jun 20 03:57:06 UTC 2024.
# Updates Price (D) / Volume(1h) (E) figures for most rows, and swaps the
# OKB / ONDO rows (44-45) back into rank order with their new figures.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value to a cell. When $asText is set, the cell is
# pre-formatted as Text so numeric-looking strings (e.g. "597.73") are not
# auto-converted to numbers by Excel, then the format is reset to the
# workbook default ("Normal" style) so no stray styling is left behind.
function Set-CellText($row, $col, $val, $asText = $false) {
    $cell = $ws.Cells.Item($row, $col)
    if ($asText) {
        $cell.NumberFormat = "@"
        $cell.Value = $val
        $cell.Style = "Normal"
    } else {
        $cell.Value = $val
    }
}

Set-CellText 2 4 "65.190.94"
Set-CellText 2 5 "  -0.34%  "
Set-CellText 3 4 "3.552.27"
Set-CellText 3 5 "  -0.28%  "
Set-CellText 4 5 "  +0.06%  "
Set-CellText 5 4 "597.73" $true
Set-CellText 5 5 "  -0.07%  "
Set-CellText 6 4 "132.97" $true
Set-CellText 7 4 "3.551.27"
Set-CellText 7 5 "  -0.33%  "
Set-CellText 8 5 "  +0.03%  "
Set-CellText 9 5 "  -0.61%  "
Set-CellText 10 5 "  -2.54%  "
Set-CellText 11 4 "7.10" $true
Set-CellText 11 5 "  -0.85%  "
Set-CellText 12 5 "  -1.01%  "
Set-CellText 13 4 "4.154.46"
Set-CellText 13 5 "  -0.16%  "
Set-CellText 14 4 "0.0000182" $true
Set-CellText 14 5 "  -2.82%  "
Set-CellText 15 4 "26.91" $true
Set-CellText 15 5 "  -0.51%  "
Set-CellText 16 4 "3.554.21"
Set-CellText 16 5 "  -0.11%  "
Set-CellText 17 5 "  -0.20%  "
Set-CellText 18 4 "65.295.45"
Set-CellText 18 5 "  -0.07%  "
Set-CellText 19 4 "9.93" $true
Set-CellText 19 5 "  -4.75%  "
Set-CellText 20 4 "14.35" $true
Set-CellText 20 5 "  +0.97%  "
Set-CellText 21 5 "  -0.79%  "
Set-CellText 22 4 "390.49" $true
Set-CellText 22 5 "  -1.67%  "
Set-CellText 23 4 "0.576" $true
Set-CellText 23 5 "  +1.07%  "
Set-CellText 24 4 "3.696.35"
Set-CellText 24 5 "  -0.14%  "
Set-CellText 25 4 "74.08" $true
Set-CellText 25 5 "  -0.87%  "
Set-CellText 26 5 "  +0.09%  "
Set-CellText 27 5 "  -1.21%  "
Set-CellText 28 4 "7.80" $true
Set-CellText 28 5 "  +0.38%  "
Set-CellText 29 4 "1.55" $true
Set-CellText 29 5 "  +24.90%  "
Set-CellText 30 5 "  -0.05%  "
Set-CellText 31 4 "8.54" $true
Set-CellText 31 5 "  +3.26%  "
Set-CellText 32 4 "2.28" $true
Set-CellText 32 5 "  +0.32%  "
Set-CellText 33 4 "3.555.11"
Set-CellText 33 5 "  -0.51%  "
Set-CellText 34 4 "24.06" $true
Set-CellText 34 5 "  -0.14%  "
Set-CellText 36 5 "  -0.50%  "
Set-CellText 37 4 "170.29" $true
Set-CellText 37 5 "  +1.01%  "
Set-CellText 38 4 "6.92" $true
Set-CellText 38 5 "  -1.75%  "
Set-CellText 39 5 "  -0.35%  "
Set-CellText 40 5 "  +0.99%  "
Set-CellText 41 4 "0.0809" $true
Set-CellText 41 5 "  +0.59%  "
Set-CellText 42 4 "0.826" $true
Set-CellText 42 5 "  -0.34%  "
Set-CellText 43 4 "26.48" $true
Set-CellText 43 5 "  -0.01%  "
Set-CellText 44 2 "OKB"
Set-CellText 44 3 "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
Set-CellText 44 4 "43.06" $true
Set-CellText 44 5 "  +0.16%  "
Set-CellText 45 2 "ONDO"
Set-CellText 45 3 "https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"
Set-CellText 45 4 "1.25" $true
Set-CellText 45 5 "  +4.99%  "
Set-CellText 46 5 "  +0.01%  "
Set-CellText 47 5 "  +0.03%  "
Set-CellText 48 5 "  -2.03%  "
Set-CellText 49 4 "2.459.07"
Set-CellText 49 5 "  +2.20%  "
Set-CellText 50 4 "6.90" $true
Set-CellText 50 5 "  +1.19%  "
Set-CellText 51 5 "  +0.75%  "
